# Add a new task row (row 10) to the "Report" sheet, mirroring the
# existing rows: Date (A) | Task Name (B) | Status (C) | Person (D)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Date value - store as the plain serial number (no time component),
# matching how the other date cells in column A are stored.
$ws.Range("A10").Value = 45762

# Copy the date's number format/alignment from the row above so the new
# cell reuses the existing date style instead of creating a new one.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B10").Value = "Added exporting Avg Volume Functions to pickle file"
$ws.Range("C10").Value = "In Progress"
$ws.Range("D10").Value = "Caleb Viverito"

# Move the active selection to D11, just below the newly added row.
$ws.Range("D11").Select()
